$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.867.52'
$ws.Range("E2").Value = '  -0.06%  '
$ws.Range("D3").Value = '3.533.36'
$ws.Range("E3").Value = '  -0.95%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '613.45'
$ws.Range("E5").Value = '  -0.10%  '
$ws.Range("D6").Value = '173.87'
$ws.Range("E6").Value = '  +0.96%  '
$ws.Range("D7").Value = '3.526.96'
$ws.Range("E7").Value = '  -1.00%  '
$ws.Range("E8").Value = '  -1.40%  '
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("E10").Value = '  -0.23%  '
$ws.Range("D11").Value = '7.41'
$ws.Range("E11").Value = '  +0.51%  '
$ws.Range("E12").Value = '  +0.13%  '
$ws.Range("D13").Value = '46.64'
$ws.Range("E13").Value = '  -0.28%  '
$ws.Range("D14").Value = '0.0000276'
$ws.Range("E14").Value = '  -0.60%  '
$ws.Range("D15").Value = '4.104.51'
$ws.Range("E15").Value = '  -0.78%  '
$ws.Range("E16").Value = '  +0.39%  '
$ws.Range("D17").Value = '615.60'
$ws.Range("E17").Value = '  -0.60%  '
$ws.Range("D18").Value = '3.532.97'
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("D19").Value = '70.883.69'
$ws.Range("E20").Value = '  +1.40%  '
$ws.Range("D21").Value = '17.78'
$ws.Range("E21").Value = '  +2.14%  '
$ws.Range("E22").Value = '  +0.37%  '
$ws.Range("D23").Value = '9.00'
$ws.Range("E23").Value = '  -5.13%  '
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").Value = '98.06'
$ws.Range("E25").Value = '  +1.10%  '
$ws.Range("E26").Value = '  -1.41%  '
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("E28").Value = '  -0.41%  '
$ws.Range("D29").Value = '33.87'
$ws.Range("E29").Value = '  +0.62%  '
$ws.Range("D30").Value = '9.15'
$ws.Range("E30").Value = '  +0.30%  '
$ws.Range("E31").Value = '  -1.15%  '
$ws.Range("D32").Value = '8.17'
$ws.Range("E32").Value = '  -4.34%  '
$ws.Range("E33").Value = '  -0.43%  '
$ws.Range("D34").Value = '6.87'
$ws.Range("E34").Value = '  -1.36%  '
$ws.Range("D35").Value = '612.30'
$ws.Range("E35").Value = '  +6.98%  '
$ws.Range("E36").Value = '  -0.66%  '
$ws.Range("D37").Value = '10.86'
$ws.Range("E37").Value = '  -0.13%  '
$ws.Range("D38").Value = '3.56'
$ws.Range("E38").Value = '  -2.12%  '
$ws.Range("E39").Value = '  +0.38%  '
$ws.Range("D40").Value = '57.04'
$ws.Range("E40").Value = '  -0.80%  '
$ws.Range("E41").Value = '  +0.18%  '
$ws.Range("D42").Value = '0.145'
$ws.Range("E42").Value = '  +0.88%  '
$ws.Range("D43").Value = '3.375.92'
$ws.Range("E43").Value = '  -0.03%  '
$ws.Range("D44").Value = '0.0₃0741'
$ws.Range("E44").Value = '  +4.96%  '
$ws.Range("E45").Value = '  -2.19%  '
$ws.Range("E46").Value = '  -2.19%  '
$ws.Range("D47").Value = '32.40'
$ws.Range("E47").Value = '  -1.96%  '
$ws.Range("E48").Value = '  -1.93%  '
$ws.Range("E49").Value = '  +0.14%  '
$ws.Range("D50").Value = '134.12'
$ws.Range("E50").Value = '  +0.30%  '
$ws.Range("E51").Value = '  -0.02%  '
